# Generate Report for Archive
# - Status text "Ready for handoff" -> "In Translation" on every sheet that
#   surfaces it (Overview!E2/F2, zh-cn!C2, de-de!C2).
# - The Status column narrows (best-fit/auto-fit) now that the replacement
#   text is shorter than the original, so the three "Status"-ish columns
#   shrink to match.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus

# Re-fit the (now narrower) status columns, mirroring what Excel does
# automatically when the cell content driving the best-fit width shrinks.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
